$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Data: YCbCr 135-degree GLCM feature grid -----------------------------
$values = @(
    @(0.030391025505013408, 0.85741118696559659, 0.76684403525661915, 0.98480457766045593),
    @(0.00002739826363958949, 0.091652967224638976, 0.99994300659013979, 0.9999863008681803),
    @(0.0050882335022645331, 0.72695331844170352, 0.97516621638259771, 0.9974558832488678),
    @(0.078717678082688841, 0.82356282891192645, 0.48493831029462359, 0.96064116095865548)
)

for ($r = 0; $r -lt 4; $r++) {
    for ($c = 0; $c -lt 4; $c++) {
        $ws.Cells.Item($r + 1, $c + 1).Value = $values[$r][$c]
    }
}

# --- Column widths ----------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 15.37890625
$ws.Columns.Item(2).ColumnWidth = 13.7109375
$ws.Columns.Item(3).ColumnWidth = 12.7109375
$ws.Columns.Item(4).ColumnWidth = 12.7109375


